$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently wraps "jth".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the "optVali" run and split it into "optVal" + "i", coloring
#    just the trailing "i" green (00B050), matching the other optimized-
#    result highlights already used elsewhere in the document.
$rng = $d.Content
$found = $rng.Find.Execute("optVali")
if ($found) {
    $start = $rng.Start
    $end = $rng.End

    $optValRange = $d.Range($start, $end - 1)
    $iRange = $d.Range($end - 1, $end)

    $iRange.Font.Color = 5287936

    # 3. Re-create the "_GoBack" bookmark around the "optVal" portion.
    $d.Bookmarks.Add("_GoBack", $optValRange)
}
